$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.224.04"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.72%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.852.52"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.25%  "

$ws.Range("E4").Value = "  -0.47%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.72%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4640"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.18%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07287"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.88%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8866"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.03%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.13"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.69%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07863"
$ws.Range("D12").Style = "Normal"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.796.92"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.11%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.392"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.97%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.508"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.60%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.07"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.47%  "

$ws.Range("E17").Value = "  -0.44%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008909"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.55%  "

$ws.Range("E19").Value = "  -0.30%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.68"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.73%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.250.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.73%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.079"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.60%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.51"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.29%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.100.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.19%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.949"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.40%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.24"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.00%  "

$ws.Range("E27").Value = "  -0.61%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.044"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.02%  "

$ws.Range("E29").Value = "  +0.20%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.033"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.81%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08824"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.79%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.142"
$ws.Range("D32").Style = "Normal"

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7669"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.22%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.166"
$ws.Range("D34").Style = "Normal"

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.523"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.67%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.718"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +10.08%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.110"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.14%  "

$ws.Range("E38").Value = "  -0.66%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05223"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.11%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.937"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.54%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.037"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.26%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5127"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.11%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1627"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.12%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.470"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.51%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4795"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.06%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.37"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.44%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.001"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.41%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "102.76"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.36%  "

$ws.Range("E49").Value = "  +0.75%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06203"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.08%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "65.48"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.05%  "
